# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-holding detail, same layout as the
# existing "2021-Q3" / "2021-Q4" sheets) positioned right before the "总计"
# (totals) sheet, and adds a corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 0: remove the original "总计" sheet so its sheetId (3) is freed up.
# We'll rebuild it at the very end (by copying the new "2022-Q1" sheet, which
# preserves all of its sheetPr/format plumbing) so the final sheetId sequence
# matches the target: 2021-Q3=1, 2021-Q4=2, 2022-Q1=3, 总计=4.
# ---------------------------------------------------------------------------
$originalTotal = $wb.Worksheets.Item(3)   # "总计"
$originalTotal.Delete()

# ---------------------------------------------------------------------------
# Step 1: create "2022-Q1" by copying "2021-Q4" (this keeps the bold/border
# header style + index-column style (cellXf #2) intact with zero extra
# styling work).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(2)             # "2021-Q4"
$ws4.Copy($null, $ws4)                    # new copy placed right after it
$newWs = $wb.Worksheets.Item(3)
$newWs.Name = "2022-Q1"

# The source sheet has 7 data rows; the target only needs 3 -> drop the rest.
$newWs.Range("A5:H8").Clear()

# Row 2 - 诺安新兴产业混合
$newWs.Range("B2").Value = "'008328"
$newWs.Range("C2").Value = "诺安新兴产业混合"
$newWs.Range("D2").Value = "'4.83"
$newWs.Range("E2").Value = "'68.49"
$newWs.Range("F2").Value = "'2.50"
$newWs.Range("G2").Value = "'0.1208"
$newWs.Range("H2").Value = 9

# Row 3 - 汇丰晋信2026周期混合
$newWs.Range("B3").Value = "'540004"
$newWs.Range("C3").Value = "汇丰晋信2026周期混合"
$newWs.Range("D3").Value = "'1.14"
$newWs.Range("E3").Value = "'31.29"
$newWs.Range("F3").Value = "'2.22"
$newWs.Range("G3").Value = "'0.0253"
$newWs.Range("H3").Value = 4

# Row 4 - 诺安改革趋势灵活配置混合
$newWs.Range("B4").Value = "'001780"
$newWs.Range("C4").Value = "诺安改革趋势灵活配置混合"
$newWs.Range("D4").Value = "'0.46"
$newWs.Range("E4").Value = "'68.34"
$newWs.Range("F4").Value = "'4.11"
$newWs.Range("G4").Value = "'0.0189"
$newWs.Range("H4").Value = 7

# ---------------------------------------------------------------------------
# Step 2: rebuild "总计" right after "2022-Q1" by copying it (again, to keep
# sheetPr / sheetFormatPr / the shared header+index-column style), then trim
# it down to the 3-column "日期 / 持有数量(只) / 持有市值(亿元)" layout.
# ---------------------------------------------------------------------------
$newWs.Copy($null, $newWs)
$totalWs = $wb.Worksheets.Item(4)
$totalWs.Name = "总计"

$totalWs.Range("E1:H8").Clear()
$totalWs.Range("A5:D8").Clear()

$totalWs.Range("B1").Value = "日期"
$totalWs.Range("C1").Value = "持有数量(只)"
$totalWs.Range("D1").Value = "持有市值(亿元)"

# New row for the 2022-Q1 summary
$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 3
$totalWs.Range("D2").Value = 0.16

# Existing rows shift down and get renumbered
$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2021-Q4"
$totalWs.Range("C3").Value = 7
$totalWs.Range("D3").Value = 0.4

$totalWs.Range("A4").Value = 2
$totalWs.Range("B4").Value = "2021-Q3"
$totalWs.Range("C4").Value = 5
$totalWs.Range("D4").Value = 0.06
